$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh cryptocurrency Price (D) and Volume(1h) (E) columns with latest feed data.
# Columns are stored as text (prices use locale-style "." thousands separators and
# percentages carry literal "%" signs with padding), so for any Price value that would
# otherwise be auto-parsed as a plain number, force the cell to Text format first so the
# original text representation (e.g. "594.40") is preserved verbatim instead of becoming 594.4.
$ws.Range("D2").Value = "65.372.09"
$ws.Range("E2").Value = "  -1.46%  "
$ws.Range("D3").Value = "3.427.31"
$ws.Range("E3").Value = "  -4.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.40"
$ws.Range("E5").Value = "  -2.03%  "
$ws.Range("E6").Value = "  -9.50%  "
$ws.Range("D7").Value = "3.425.60"
$ws.Range("E7").Value = "  -4.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.489"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.47"
$ws.Range("E10").Value = "  -4.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.121"
$ws.Range("E11").Value = "  -10.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.377"
$ws.Range("E12").Value = "  -8.76%  "
$ws.Range("D13").Value = "4.003.25"
$ws.Range("E13").Value = "  -4.61%  "
$ws.Range("E14").Value = "  -12.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.35"
$ws.Range("E15").Value = "  -10.62%  "
$ws.Range("D16").Value = "65.337.51"
$ws.Range("E16").Value = "  -1.61%  "
$ws.Range("D17").Value = "3.415.04"
$ws.Range("E17").Value = "  -4.59%  "
$ws.Range("E18").Value = "  -2.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.85"
$ws.Range("E19").Value = "  -10.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.72"
$ws.Range("E20").Value = "  -9.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.67"
$ws.Range("E21").Value = "  -7.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "390.95"
$ws.Range("E22").Value = "  -7.51%  "
$ws.Range("E23").Value = "  -6.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.542"
$ws.Range("E24").Value = "  -11.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").Value = "3.568.29"
$ws.Range("E26").Value = "  -4.40%  "
$ws.Range("E27").Value = "  -12.57%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("E29").Value = "  -9.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.14"
$ws.Range("E30").Value = "  -12.85%  "
$ws.Range("E31").Value = "  -12.73%  "
$ws.Range("D32").Value = "3.431.80"
$ws.Range("E32").Value = "  -4.27%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  -7.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "22.63"
$ws.Range("E35").Value = "  -9.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "172.55"
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("E37").Value = "  -14.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.80"
$ws.Range("E38").Value = "  -12.10%  "
$ws.Range("E39").Value = "  -9.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.81"
$ws.Range("E40").Value = "  -13.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0767"
$ws.Range("E41").Value = "  -9.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.810"
$ws.Range("E42").Value = "  -8.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "43.40"
$ws.Range("E43").Value = "  -5.34%  "
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.38"
$ws.Range("E45").Value = "  -15.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.61"
$ws.Range("E46").Value = "  -13.38%  "
$ws.Range("E47").Value = "  -3.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.94"
$ws.Range("E48").Value = "  -6.59%  "
$ws.Range("E49").Value = "  -8.95%  "
$ws.Range("E50").Value = "  -15.61%  "
$ws.Range("D51").Value = "2.183.64"
